$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.417.84"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "1.840.02"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "261.87"
$ws.Range("E5").Value = "  -5.60%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.5192"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").Value = "0.3263"
$ws.Range("E8").Value = "  -4.62%  "
$fmt_D9 = $ws.Range("D9").NumberFormat
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06780"
$ws.Range("D9").NumberFormat = $fmt_D9
$ws.Range("D10").Value = "18.68"
$ws.Range("E10").Value = "  -6.75%  "
$ws.Range("D11").Value = "0.7731"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").Value = "0.07739"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "1.862.49"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "87.97"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "5.001"
$ws.Range("E15").Value = "  -3.48%  "
$fmt_D16 = $ws.Range("D16").NumberFormat
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").NumberFormat = $fmt_D16
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "13.93"
$ws.Range("E17").Value = "  -4.38%  "
$fmt_D18 = $ws.Range("D18").NumberFormat
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").NumberFormat = $fmt_D18
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.000007935"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").Value = "26.453.36"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "2.072.88"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").Value = "4.601"
$ws.Range("D23").Value = "9.528"
$ws.Range("E23").Value = "  -5.04%  "
$ws.Range("D24").Value = "6.008"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "145.15"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").Value = "2.184"
$ws.Range("E26").Value = "  -8.18%  "
$ws.Range("D27").Value = "1.654"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "16.94"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").Value = "111.76"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").Value = "4.191"
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("D31").Value = "4.126"
$ws.Range("E31").Value = "  -4.66%  "
$fmt_D32 = $ws.Range("D32").NumberFormat
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08700"
$ws.Range("D32").NumberFormat = $fmt_D32
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "0.04809"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").Value = "0.7164"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "3.084"
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("D38").Value = "0.01779"
$ws.Range("E38").Value = "  -4.18%  "
$fmt_D39 = $ws.Range("D39").NumberFormat
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.220"
$ws.Range("D39").NumberFormat = $fmt_D39
$ws.Range("E39").Value = "  -4.40%  "
$ws.Range("D40").Value = "0.4834"
$ws.Range("E40").Value = "  -6.08%  "
$ws.Range("D41").Value = "112.27"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Value = "0.9004"
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").Value = "6.074"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "7.719"
$ws.Range("E45").Value = "  -4.85%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.05914"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4145"
$ws.Range("E47").Value = "  -7.24%  "
$ws.Range("D48").Value = "9.028"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").Value = "35.02"
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").Value = "0.1218"
$ws.Range("E50").Value = "  -9.19%  "
$ws.Range("D51").Value = "0.8859"
$ws.Range("E51").Value = "  +0.01%  "
